$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Online")
$ws.Range("A3:H22").ClearContents()
